# Apply crypto price/volume updates scraped on Tue Nov 28 00:29:28 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes plain text that must survive as text even when it
# looks numeric (e.g. "1.00"), without leaving the cell restyled -
# the temporary "@" text format is reset back to the default "Normal"
# style right after the value is written.
function Set-CellText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "37.275.87"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "2.031.72"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  -0.37%  "
Set-CellText "D5" "227.88"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("E7").Value = "  +0.05%  "
Set-CellText "D8" "55.37"
$ws.Range("E8").Value = "  -4.23%  "
$ws.Range("E9").Value = "  -2.60%  "
Set-CellText "D10" "0.0796"
$ws.Range("E10").Value = "  +0.27%  "
Set-CellText "D11" "0.102"
$ws.Range("E11").Value = "  -5.35%  "
$ws.Range("D12").Value = "2.334.27"
$ws.Range("E12").Value = "  -1.54%  "
Set-CellText "D13" "14.29"
$ws.Range("E13").Value = "  -4.60%  "
Set-CellText "D14" "20.41"
$ws.Range("E14").Value = "  -3.90%  "
Set-CellText "D15" "0.746"
$ws.Range("E15").Value = "  -2.62%  "
$ws.Range("E16").Value = "  -3.44%  "
$ws.Range("D17").Value = "2.025.97"
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").Value = "37.260.74"
$ws.Range("E18").Value = "  -0.54%  "
Set-CellText "D19" "6.27"
$ws.Range("E19").Value = "  +0.05%  "
Set-CellText "D20" "69.27"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("D21").Value = "0.0₃0824"
$ws.Range("E21").Value = "  -1.33%  "
Set-CellText "D22" "224.57"
$ws.Range("E22").Value = "  -1.79%  "
Set-CellText "D23" "1.00"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("E25").Value = "  -6.05%  "
Set-CellText "D26" "9.33"
$ws.Range("E26").Value = "  -5.73%  "
Set-CellText "D27" "165.49"
$ws.Range("E27").Value = "  -2.27%  "
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-CellText "D29" "18.84"
$ws.Range("E29").Value = "  -2.52%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-CellText "D30" "1.36"
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("E31").Value = "  -2.75%  "
Set-CellText "D32" "4.56"
$ws.Range("E32").Value = "  -1.38%  "
Set-CellText "D33" "0.0617"
$ws.Range("E33").Value = "  -2.27%  "
Set-CellText "D34" "4.49"
$ws.Range("E34").Value = "  -3.75%  "
$ws.Range("E35").Value = "  -5.99%  "
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  -4.92%  "
Set-CellText "D39" "5.51"
$ws.Range("E39").Value = "  +4.29%  "
$ws.Range("E40").Value = "  -3.57%  "
$ws.Range("D41").Value = "1.479.58"
$ws.Range("E41").Value = "  -1.89%  "
Set-CellText "D42" "97.00"
$ws.Range("E42").Value = "  -1.33%  "
Set-CellText "D43" "16.82"
$ws.Range("E43").Value = "  -2.28%  "
Set-CellText "D44" "0.0928"
$ws.Range("E44").Value = "  -2.69%  "
Set-CellText "D45" "2.79"
$ws.Range("E45").Value = "  -4.39%  "
$ws.Range("E46").Value = "  -5.28%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-CellText "D47" "7.27"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-CellText "D48" "1.02"
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").Value = "2.219.04"
$ws.Range("E50").Value = "  -1.62%  "
Set-CellText "D51" "3.58"
$ws.Range("E51").Value = "  -11.44%  "

Write-Output "Applied 92 cell updates"
